$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "43.852.11"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.225.85"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.58"
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.24"
$ws.Range("E6").Value = "  +11.84%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.31"
$ws.Range("E10").Value = "  +9.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.53"
$ws.Range("E12").Value = "  +7.85%  "
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "2.559.15"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.81"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "2.217.42"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "43.816.41"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.05"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("E22").Value = "  +4.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.46"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("E25").Value = "  +19.59%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.54"
$ws.Range("E28").Value = "  +6.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.23"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.02"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0907"
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.59"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.34"
$ws.Range("E39").Value = "  +18.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.51"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.66"
$ws.Range("E41").Value = "  +9.25%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.31"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0993"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.31"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.31"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("E48").Value = "  +7.48%  "
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.441"
$ws.Range("E50").Value = "  -9.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
